$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new data row (row 47) with the latest metric reading
$ws.Cells.Item(47, 1).Value = "2025-04-29 05:14:51"
$ws.Cells.Item(47, 2).Value = 151
